$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4666
$ws.Range("J13").Value = 4666
$ws.Range("L13").Value = 4666
$ws.Range("N13").Value = -5004

$ws.Range("H92").Value = 650.1667
$ws.Range("I92").Value = 655.3
$ws.Range("K92").Value = 655.3
$ws.Range("M92").Value = 592.7

$ws.Range("H94").Value = 10337.5
$ws.Range("I94").Value = 10337.5
$ws.Range("K94").Value = 10337.5
$ws.Range("M94").Value = -9886.5

$ws.Range("H106").Value = 3966.3333
$ws.Range("I106").Value = 3966.3333
$ws.Range("K106").Value = 3966.3333
$ws.Range("M106").Value = -3335.3333

$ws.Range("H138").Value = 1475.2979
$ws.Range("I138").Value = 889.931
$ws.Range("J138").Value = 2418.389
$ws.Range("K138").Value = 2669.793
$ws.Range("L138").Value = 7255.167
$ws.Range("M138").Value = 2470.207
$ws.Range("N138").Value = -17535.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2108.3
$ws.Range("I61").Value = 1763.8334
$ws.Range("K61").Value = 1763.8334
$ws.Range("M61").Value = -1551.8334

$ws.Range("H74").Value = 2149.182
$ws.Range("I74").Value = 1383.3125
$ws.Range("K74").Value = 1383.3125
$ws.Range("M74").Value = -509.3125

$ws.Range("H77").Value = 2149.182
$ws.Range("I77").Value = 1383.3125
$ws.Range("K77").Value = 6916.5625
$ws.Range("M77").Value = -2548.5625

$ws.Range("H102").Value = 184626.83
$ws.Range("I102").Value = 334500
$ws.Range("J102").Value = 34753.668
$ws.Range("K102").Value = 334500
$ws.Range("L102").Value = 34753.668
$ws.Range("M102").Value = -332878
$ws.Range("N102").Value = -37997.668

$ws.Range("H104").Value = 28873.428
$ws.Range("J104").Value = 28873.428
$ws.Range("L104").Value = 28873.428
$ws.Range("N104").Value = -35861.428

$ws.Range("H108").Value = 88996.8
$ws.Range("J108").Value = 88996.8
$ws.Range("L108").Value = 88996.8
$ws.Range("N108").Value = -96676.8

$ws.Range("H121").Value = 51655.555
$ws.Range("J121").Value = 51655.555
$ws.Range("L121").Value = 51655.555
$ws.Range("N121").Value = -55149.555

$ws.Range("H132").Value = 1367.7084
$ws.Range("I132").Value = 1255.6666
$ws.Range("K132").Value = 3766.9998
$ws.Range("M132").Value = -1236.9998

$ws.Range("H136").Value = 2108.3
$ws.Range("I136").Value = 1763.8334
$ws.Range("K136").Value = 5291.5002
$ws.Range("M136").Value = -2741.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 30496
$ws.Range("J2").Value = 30496
$ws.Range("L2").Value = 30496
$ws.Range("N2").Value = -30722

$ws.Range("H50").Value = 58676
$ws.Range("J50").Value = 60667.75
$ws.Range("L50").Value = 60667.75
$ws.Range("N50").Value = -61815.75

$ws.Range("H55").Value = 35897.4
$ws.Range("J55").Value = 35897.4
$ws.Range("L55").Value = 35897.4
$ws.Range("N55").Value = -36443.4

$ws.Range("H105").Value = 52855.4
$ws.Range("I105").Value = 101650.8
$ws.Range("J105").Value = 4060
$ws.Range("K105").Value = 101650.8
$ws.Range("L105").Value = 4060
$ws.Range("M105").Value = -99903.8
$ws.Range("N105").Value = -7554

$ws.Range("H115").Value = 74425.86
$ws.Range("J115").Value = 76496.664
$ws.Range("L115").Value = 76496.664
$ws.Range("N115").Value = -79630.664

$ws.Range("H141").Value = 80555.55499999999
$ws.Range("J141").Value = 80555.55499999999
$ws.Range("L141").Value = 80555.55499999999
$ws.Range("N141").Value = -90915.55499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 683.3333
$ws.Range("I2").Value = 683.3333
$ws.Range("K2").Value = 683.3333
$ws.Range("M2").Value = -570.3333

$ws.Range("H3").Value = 6666.6665
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10226
$ws.Range("M3").Value = -4887

$ws.Range("H18").Value = 28890.334
$ws.Range("J18").Value = 28890.334
$ws.Range("L18").Value = 28890.334
$ws.Range("N18").Value = -29350.334

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H116").Value = 85486.25
$ws.Range("J116").Value = 85486.25
$ws.Range("L116").Value = 85486.25
$ws.Range("N116").Value = -94664.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 28.166666
$ws.Range("I2").Value = 18.181818
$ws.Range("J2").Value = 43.857143
$ws.Range("K2").Value = 109.090908
$ws.Range("L2").Value = 263.142858
$ws.Range("M2").Value = 3.909092000000001
$ws.Range("N2").Value = -489.142858

$ws.Range("H111").Value = 481.25
$ws.Range("I111").Value = 481.25
$ws.Range("K111").Value = 1443.75
$ws.Range("M111").Value = 1623.25

$ws.Range("H140").Value = 1091.7273
$ws.Range("I140").Value = 850.9
$ws.Range("K140").Value = 2552.7
$ws.Range("M140").Value = 2627.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 52120.6
$ws.Range("J103").Value = 50151
$ws.Range("L103").Value = 50151
$ws.Range("N103").Value = -52495

$ws.Range("H114").Value = 60094.832
$ws.Range("J114").Value = 60094.832
$ws.Range("L114").Value = 60094.832
$ws.Range("N114").Value = -68772.83199999999

$ws.Range("H135").Value = 57142.285
$ws.Range("J135").Value = 57142.285
$ws.Range("L135").Value = 57142.285
$ws.Range("N135").Value = -67282.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1316.8889

$ws.Range("H27").Value = 1316.8889

$ws.Range("H46").Value = 3428.1155
$ws.Range("I46").Value = 1906.25
$ws.Range("J46").Value = 3704.818
$ws.Range("K46").Value = 1906.25
$ws.Range("L46").Value = 3704.818
$ws.Range("M46").Value = -1718.25
$ws.Range("N46").Value = -4080.818

$ws.Range("H117").Value = 64323
$ws.Range("J117").Value = 69097.336
$ws.Range("L117").Value = 69097.336
$ws.Range("N117").Value = -78275.336

$ws.Range("H123").Value = 74797.25
$ws.Range("J123").Value = 78284
$ws.Range("L123").Value = 78284
$ws.Range("N123").Value = -88084

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 37435.5
$ws.Range("J121").Value = 37435.5
$ws.Range("L121").Value = 37435.5
$ws.Range("N121").Value = -40929.5
